$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @('John', 25, 'New York', '25th street'),
    @('Alice', 30, 'Los Angeles', 'apartment 25, 5th floor'),
    @('Bob', 22, 'Chicago', '26th street'),
    @('Charlie', 28, 'Houston', 'apartment 25, 5th floor'),
    @('David', 35, 'Phoenix', '27th street'),
    @('Emma', 40, 'Philadelphia', 'apartment 25, 5th floor'),
    @('Fathi', 19, 'San Antonio', '28th street'),
    @('Grace', 21, 'San Diego', 'apartment 25, 5th floor'),
    @('Henry', 45, 'Dallas', '29th street'),
    @('Ivy', 50, 'San Jose', 'apartment 25, 5th floor'),
    @('Jack', 33, 'Austin', '30th street'),
    @('Karen', 27, 'Jacksonville', 'apartment 25, 5th floor'),
    @('Leo', 31, 'Fort Worth', '31st street'),
    @('Mona', 29, 'Columbus', 'apartment 25, 5th floor'),
    @('Nathan', 26, 'Charlotte', '32nd street'),
    @('Olivia', 38, 'Indianapolis', 'apartment 25, 5th floor'),
    @('Paul', 41, 'Seattle', '33rd street'),
    @('Quinn', 20, 'Denver', 'apartment 25, 5th floor'),
    @('Rachel', 32, 'Washington', '34th street'),
    @('Steve', 37, 'Boston', 'apartment 25, 5th floor')
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 1
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}
